# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append week-16 yardage samples to the four running play-logs.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "5 -2 4 2 -1 1 5 2 3 7 3 5 2 4 2 6 1 4 5 2 15 -2 2 6 41 6 2 7 5 3 0 9 -1 -6 2 -9 7 2 2 5 11 1 1 2 8 7 2 5 12 -3 15 2 -1 3 4 0 4 6 -3 2 4 2 6 11 4 4 1 1 6 -1 0 4 0 1 7 8 1 17 20 1 1 -3 2 4 1 4 2 1 6 6 4 6 5 1 4 3 4 1 14 2 4 1 0 0 -1 4 0 2 4 -2 7 3 4 4 14 9 2 3 0 2 2 0 6 -1 2 4 2 2 5 4 14 -1 7 2 2 0 1 9 3 0 -2 -2 7 11 3 4 2 3 6 2 2 -2 4 19 0 4 9 2 0 2 -3 11 2 2 6 2 14 5 4 0 3 5 4 10 2 2 6 4 16 5 5 6 2 -2 2 6 1 6 20 5 2 3 5 2 0 5 7 5 2 3 8 -3 4 11 12 2 2 8 6 0 1 3 16 1 2 12 0 9 6 5 0 2 14 2 2 1 1 32 6 -1 -3 7 -4 3 4 13 -1 1 4 3 5 0 1 -1 3 3 10 -2 0 1 0 7 0 23 17 1 0 9 13 3 7 8 4 4 2 0 7 -1 2 7 5 1 7 5 7 4 9 5 -3 13 8 5 8 6 8 7 1 6 5 4 5 5 28 4 5 0 5 1 -3 0 6 6 3 5 1 3 2 0 1 5 5 3 4 2 1 6 6 -2 -1 -4 10 3 7 1 6 5 0 5 0 2 2 4 4 2 7"
$ws.Range("C2").Value = "-1 2 5 2 9 2 2 -2 8 7 3 0 5 4 6 6 2 2 12 3 8 0 2 4 70 3 -1 2 1 6 2 5 0 5 1 10 2 2 2 2 7 14 3 1 3 1 11 3 4 1 0 5 2 11 1 3 11 5 10 3 3 3 2 0 2 1 3 -1 6 10 6 3 3 3 2 5 4 2 3 8 1 3 2 5 8 1 8 5 3 4 6 9 2 9 4 4 4 8 7 1 0 3 2 8 1 5 2 15 4 4 6 6 17 4 5 3 6 6 3 5 12 0 0 14 0 4 3 7 -2 5 3 14 0 10 2 2 7 4 13 2 2 2 3 2 1 0 22 9 15 0 4 3 4 2 2 8 0 8 6 10 5 3 1 3 0 1 9 6 5 2 4 5 0 3 3 0 3 6 13 7 -3 2 -2 2 2 4 7 -4 0 1 3 2 5 2 10 -3 4 2 4 4 9 7 5 9 3 8 12 -1 3 4 6 0 -1 1 8 6 -1 11 6 13 3 0 13 3 7 2 1 14 4 1 0 3 4 21 2 3 7 10 7 6 3 2 5 10 3 3 4 6 2 7 -1 3 4 9 2 2 1 3 -3 2 -1 -1 1 7 4 5 3 2 11 5 7 2 14 18 6 1 27 3 5 12 5 4 13 5 13 12 1 1 6 3 4 3 4 0 2 2 4 11 2 0 6 1 6 2 3 0 2 3 4 1 3 4 4 2 3 2 2 0 14 9 5 5 9 1 1 0 2 17 2 3 3 12 4 -1 6 4 8 2 4 2 6 1 3 2 10 9 0 3 2 0 2 4 3 -1 7 8 4 13 17 2 3 2 0 2 3 10 7 9 12 2 5 0 0 2 1 -2 1 12 3 1 11 0 3 15 5 0 13 4 1 1 2 5 5 2 3 9 1 4 4 6 6 6 4 5 3 3 3 6 0"
$ws.Range("B3").Value = "42 -6 4 17 1 6 8 37 2 15 20 8 17 14 6 16 16 6 9 15 8 6 12 16 13 4 4 8 2 6 17 3 4 2 8 12 19 33 13 21 10 8 21 13 14 8 38 19 8 18 12 9 7 4 16 0 13 13 11 10 12 3 7 5 16 4 9 10 13 -4 6 18 4 27 51 19 7 8 17 20 8 4 5 21 1 54 8 5 11 3 28 17 8 18 1 23 38 5 9 18 28 6 35 6 2 26 5 13 13 16 14 13 13 11 7 13 17 8 3 16 2 9 4 14 7 8 3 10 10 11 6 10 3 10 1 11 12 11 11 7 6 2 15 -2 6 18 7 1 27 15 7 16 15 1 6 8 15 -4 2 14 16 2 5 8 9 12 5 50 19 3 11 1 12 41 2 5 4 22 4 6 11 6 5 7 16 4 3 3 9 5 30 6 12 7 -7 14 3 14 2 9 1 -1 6 9 12 9 4 3 5 5 0 11 21 4 28 8 11 8 -1 7 2 4 5 5 1 20 7 18 19 6 14 9 0 10 17 17 1 18 18 5 4 10 8 3 15 -3 7 -5 4 18 20 -1 5 12 16 6 11 17 8 3 9 8 2 10 11 11 5 15 6 60 3 7 -2 1 17 3 9 6 18 10 1 8 24 8 2 16 4 15 0 3 4 9 15 9 7 7 8 14 7 7 10 36 -4 8 8 8 -2 7 6 8 2 9 11 5 2 0 4 5 9 6 4 8 2"
$ws.Range("C3").Value = "15 13 7 6 1 9 5 25 -4 5 11 6 2 6 17 16 2 10 6 8 20 7 22 13 3 7 5 22 12 4 4 16 11 8 12 0 24 12 7 6 10 7 4 4 8 18 9 0 3 10 7 56 19 2 7 9 5 7 6 3 8 19 3 4 9 20 4 6 12 14 3 4 1 12 3 -1 10 26 4 6 14 4 1 1 28 25 -3 11 9 6 18 2 11 6 5 6 19 15 58 9 21 15 3 16 8 6 8 6 2 49 6 8 -2 11 21 7 15 24 11 30 15 4 8 14 18 22 35 7 30 16 6 28 10 15 9 4 4 -6 9 -4 16 25 4 4 5 18 5 6 13 4 2 9 1 3 5 5 12 21 5 12 6 13 7 5 1 7 7 19 10 3 9 12 7 0 15 5 6 6 7 11 9 17 9 3 11 12 9 18 8 4 9 10 10 15 14 24 2 6 6 20 6 2 8 3 7 8 10 20 7 2 15 9 5 9 20 24 22 11 2 3 12 3 19 9 17 8 3 10 8 15 3 13 7 10 8 6 6 6 7 4 35 7 12 16 3 13 9 17 14 14 18 4 5 8 3 21 5 3 14 0 19 2 0 13 4 21 6 10 8 7 14 11 9 5 11 1 4 12 8 6 -2 6 8 4 8 5 9 2 9 17 4 6 5 3 17 25 11 13 3 2 16 17 15 7 5 11 9 -2 15 12 10 25 7 2 7 5 17 59 7 36 5 11 8 3 1 9 11 4 2 19 4 9 9 6 5 11 3 3 6 7 15 9 3 7 3 10 16 10 1 11 7 9 7 6 2 3 23 13 46 1 9 1 6 39 4 3 16 10 5"

# ---------------------------------------------------------------------------
# OFF sheet: week-16 box score added to row 2 (Home), season sim added to row 3 (Road).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 158
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 70
$ws.Range("G2").Value = 39
$ws.Range("H2").Value = 9
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 26
$ws.Range("N2").Value = 16
$ws.Range("O2").Value = 25
$ws.Range("P2").Value = 11

$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 196
$ws.Range("E3").Value = 27
$ws.Range("F3").Value = 105
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 25
$ws.Range("I3").Value = 71
$ws.Range("J3").Value = 54
$ws.Range("L3").Value = 305
$ws.Range("M3").Value = 179
$ws.Range("Q3").Value = 504

# ---------------------------------------------------------------------------
# DEF sheet: same treatment, defensive side of the ball.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 202
$ws.Range("F2").Value = 60
$ws.Range("G2").Value = 70
$ws.Range("H2").Value = 4
$ws.Range("J2").Value = 28
$ws.Range("N2").Value = 11
$ws.Range("O2").Value = 20

$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 178
$ws.Range("E3").Value = 25
$ws.Range("F3").Value = 107
$ws.Range("G3").Value = 37
$ws.Range("H3").Value = 17
$ws.Range("I3").Value = 59
$ws.Range("J3").Value = 62
$ws.Range("L3").Value = 300
$ws.Range("M3").Value = 206
$ws.Range("Q3").Value = 559

# ---------------------------------------------------------------------------
# ST sheet: kick/return totals (row 2) plus the RM sample logs (rows 4-6 / D3-D5).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 60
$ws.Range("D2").Value = 64
$ws.Range("F2").Value = 360
$ws.Range("G2").Value = 343
$ws.Range("N2").Value = 51
$ws.Range("O2").Value = 32

$ws.Range("B4").Value = "63 62 65 64 68 56 53 68 67 58 65 47 54 62 56 52 61 65 64 62 56 64 51 52 63 61 57 64 58 71 62 60 65 51"
$ws.Range("B5").Value = "21 27 16 23 24 10 16 16 24 18 19 7 16 13 10 14 26 32 22 20 18 13 15 13 15 26 15 17 0 29 21 29 20 10"
$ws.Range("B6").Value = "28 38 20 15 33 21 30 22 23 37 22 11 23 12 19 18 18 22 17 17 15 17 21"
$ws.Range("D3").Value = "43 58 61 28 42 53 55 49 41 34 50 43 51 33 30 40 50 57 57 33 63 36 44 48 50 51 45 51 51 34 34 42 58 46 52 39 51 51 46 33 36 42 39 36 41 45 46 44 36 50 46 51 63 18 50 51 47 48 28 38 33 33 31 36"
$ws.Range("D4").Value = "0 13 11 0 0 0 15 9 0 0 0 9 0 0 0 17 0 8 28 0 11 0 10 0 0 4 0 0 0 0 0 0 0 5 0 0 9 0 4 0 2 3 9 0 0 7 0 0 0 0 0 0 0 0 0 0 -4 0 0 14 0 39 0 1"
$ws.Range("D5").Value = "0 2 11 8 0 0 0 0 0 0 19 17 0 26 0 7 0 0 0 0 0 0 2 12 6 0 4 0 0 0 12 14 0 0 0 9 0 0 0 15 11 0 -4 0 8 0 0 0 0 0 0 0 3 13 0"

# ---------------------------------------------------------------------------
# TURNS sheet: Road row (row 3) turnover totals.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("B3").Value = 10
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 8

# ---------------------------------------------------------------------------
# PEN sheet: False start penalty totals.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 18
$ws.Range("D2").Value = 5
